$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# "Generate Report for Handback" — a handback happened for both UUID files in
# both the zh-cn and de-de localization sheets. For each sheet/row we now
# know the "Latest Target File" (F) and "Latest Handback File" (G), and the
# zh-cn sheet's status flips from "Ready for handoff" to
# "Handed back: in sync with en-US" (shared text, so it updates both rows).
# Handback datetimes are recorded per-language in column H.
# ---------------------------------------------------------------------------

function Update-HandbackSheet {
    param(
        $ws,
        $lang,                  # "zh-cn" or "de-de"
        $handbackDateTime,      # text written into column H for both rows
        $statusText             # $null to leave Status (col C) untouched
    )

    $uuid1 = "2aedfd63-731e-4312-be44-564ec605202a"
    $uuid2 = "b6121ff3-8af0-4070-aa99-c314affb1e78"
    $hash1 = "afe9b6edca296b8f1d16ddb7cfcfc392d90a4a4b"
    $hash2 = "4f866adc729a3108880e502bbfa0a278148fedf1"

    $md1 = "$uuid1.md"
    $md2 = "$uuid2.md"
    $xlf1 = "$uuid1.$hash1.$lang.xlf"
    $xlf2 = "$uuid2.$hash2.$lang.xlf"

    $urlMd1 = "https://github.com/OpenLocalizationTest/oltest/blob/a9e8475976d4403b0aee7ddf8878a2fc95746668/e2e/$md1"
    $urlMd2 = "https://github.com/OpenLocalizationTest/oltest/blob/a9e8475976d4403b0aee7ddf8878a2fc95746668/e2e/$md2"

    if ($lang -eq "zh-cn") {
        $urlXlf1 = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/23cd242380876c60eef70a06a1113a5aeaf1135f/ol-handoff/OpenLocalizationTestOrg/oltest-zhcn-fly/yuwzho/ht/$xlf1"
        $urlXlf2 = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/23cd242380876c60eef70a06a1113a5aeaf1135f/ol-handoff/OpenLocalizationTestOrg/oltest-zhcn-fly/yuwzho/ht/$xlf2"
    } else {
        $urlXlf1 = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/3d27a0f1ed1651de55ae9c1c1c66e9ea32a5b117/ol-handoff/OpenLocalizationTestOrg/oltest-dede-fly/yuwzho/ht/$xlf1"
        $urlXlf2 = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/3d27a0f1ed1651de55ae9c1c1c66e9ea32a5b117/ol-handoff/OpenLocalizationTestOrg/oltest-dede-fly/yuwzho/ht/$xlf2"
    }

    if ($statusText) {
        $ws.Range("C2").Value = $statusText
        $ws.Range("C3").Value = $statusText
    }

    # New columns: F = Latest Target File, G = Latest Handback File.
    # They mirror the existing Latest Handoff File (A = md, D = xlf).
    $ws.Range("F2").Value = $md1
    $ws.Range("G2").Value = $xlf1
    $ws.Range("F3").Value = $md2
    $ws.Range("G3").Value = $xlf2

    # Latest Handback DateTime.
    $ws.Range("H2").Value = $handbackDateTime
    $ws.Range("H3").Value = $handbackDateTime

    # Rebuild hyperlinks so the new F/G links land in the expected order,
    # immediately following the column they mirror.
    $ws.Hyperlinks.Delete()
    $ws.Hyperlinks.Add($ws.Range("A2"), $urlMd1, "", "", $md1)   | Out-Null
    $ws.Hyperlinks.Add($ws.Range("D2"), $urlXlf1, "", "", $xlf1) | Out-Null
    $ws.Hyperlinks.Add($ws.Range("F2"), $urlMd1, "", "", $md1)   | Out-Null
    $ws.Hyperlinks.Add($ws.Range("G2"), $urlXlf1, "", "", $xlf1) | Out-Null
    $ws.Hyperlinks.Add($ws.Range("A3"), $urlMd2, "", "", $md2)   | Out-Null
    $ws.Hyperlinks.Add($ws.Range("D3"), $urlXlf2, "", "", $xlf2) | Out-Null
    $ws.Hyperlinks.Add($ws.Range("F3"), $urlMd2, "", "", $md2)   | Out-Null
    $ws.Hyperlinks.Add($ws.Range("G3"), $urlXlf2, "", "", $xlf2) | Out-Null
}

$wsZhCn = $wb.Worksheets.Item("zh-cn")
Update-HandbackSheet $wsZhCn "zh-cn" "2016-03-19 04:09:30" "Handed back: in sync with en-US"

$wsDeDe = $wb.Worksheets.Item("de-de")
Update-HandbackSheet $wsDeDe "de-de" "2016-03-19 04:09:44" $null
